$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5: correct the casing of the AlphaEdit model name ---
$ws.Range("A5").Value = "AlphaEdit"

# --- Append new rows (8-12) with fresh paper entries ---
# Copy the formatting of an existing data row (row 7) down onto the
# new rows first, so font/alignment/row-height match the rest of the table.
# Row 10 has no model name (column A left blank), so it is formatted
# separately over just B:D to avoid fabricating an empty A10 cell.
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A8:D9").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:D12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B7:D7").Copy() | Out-Null
$ws.Range("B10:D10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Rows("8:12").RowHeight = 18.75

# Row 8
$ws.Range("A8").Value = "LiveEdit"
$ws.Range("B8").Value = "CVPR 2025"
$ws.Range("C8").Value = "Lifelong Knowledge Editing for Vision Language Models with  Low-Rank Mixture-of-Experts"
$ws.Range("D8").Value = "多模态大模型"

# Row 9
$ws.Range("A9").Value = "MC-MKE"
$ws.Range("B9").Value = "ACL2025"
$ws.Range("C9").Value = "MC-MKE: A Fine-Grained Multimodal Knowledge Editing Benchmark  Emphasizing Modality Consistency"
$ws.Range("D9").Value = "多模态大模型benchmark"

# Row 10 (no model name given)
$ws.Range("B10").Value = "ACM MM2025"
$ws.Range("C10").Value = "Deciphering Functions of Neurons in Vision-Language Models"
$ws.Range("D10").Value = "多模态大模型神经元的机理分析"

# Row 11
$ws.Range("A11").Value = "AcE"
$ws.Range("B11").Value = "arXiv"
$ws.Range("C11").Value = "ACE: Attribution-Controlled Knowledge Editing for Multi-hop Factual Recall"
$ws.Range("D11").Value = "多模态大模型"

# Row 12
$ws.Range("A12").Value = "STEAM"
$ws.Range("B12").Value = "arXiv"
$ws.Range("C12").Value = "STEAM: A Semantic-Level Knowledge Editing Framework for Large Language Models"
$ws.Range("D12").Value = "多模态大模型"

# --- Widen column D so the longer entries remain fully visible ---
$ws.Columns("D").ColumnWidth = 37

# --- Restore the on-screen selection left by the author ---
$ws.Range("C27").Select() | Out-Null
